$d = $word.ActiveDocument
# --- paragraph 2 ---
$p = $d.Paragraphs.Item(2)
$oldLen = $p.Range.End - $p.Range.Start - 1
$insPoint = $p.Range.Start
$ins = $d.Range($insPoint, $insPoint)
$xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>
<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">
<pkg:part pkg:name=`"/_rels/.rels`" pkg:contentType=`"application/vnd.openxmlformats-package.relationships+xml`" pkg:padding=`"512`">
<pkg:xmlData>
<Relationships xmlns=`"http://schemas.openxmlformats.org/package/2006/relationships`">
<Relationship Id=`"rId1`" Type=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument`" Target=`"word/document.xml`"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">
<pkg:xmlData>
<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">
<w:body>
<w:p><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`">docker </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`">images  </w:t></w:r><w:r><w:t>=</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t>&gt; show all images</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>"
$ins.InsertXML($xml)
$p2 = $d.Paragraphs.Item(2)
$newEnd = $p2.Range.End - 1
$delStart = $newEnd - $oldLen
$delRng = $d.Range($delStart, $newEnd)
$delRng.Delete()
# --- paragraph 3 ---
$p = $d.Paragraphs.Item(3)
$oldLen = $p.Range.End - $p.Range.Start - 1
$insPoint = $p.Range.Start
$ins = $d.Range($insPoint, $insPoint)
$xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>
<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">
<pkg:part pkg:name=`"/_rels/.rels`" pkg:contentType=`"application/vnd.openxmlformats-package.relationships+xml`" pkg:padding=`"512`">
<pkg:xmlData>
<Relationships xmlns=`"http://schemas.openxmlformats.org/package/2006/relationships`">
<Relationship Id=`"rId1`" Type=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument`" Target=`"word/document.xml`"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">
<pkg:xmlData>
<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">
<w:body>
<w:p><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`">docker </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t>ps</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`"> -a </w:t></w:r><w:r><w:t>=&gt; all containers</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>"
$ins.InsertXML($xml)
$p2 = $d.Paragraphs.Item(3)
$newEnd = $p2.Range.End - 1
$delStart = $newEnd - $oldLen
$delRng = $d.Range($delStart, $newEnd)
$delRng.Delete()
# --- paragraph 4 ---
$p = $d.Paragraphs.Item(4)
$oldLen = $p.Range.End - $p.Range.Start - 1
$insPoint = $p.Range.Start
$ins = $d.Range($insPoint, $insPoint)
$xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>
<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">
<pkg:part pkg:name=`"/_rels/.rels`" pkg:contentType=`"application/vnd.openxmlformats-package.relationships+xml`" pkg:padding=`"512`">
<pkg:xmlData>
<Relationships xmlns=`"http://schemas.openxmlformats.org/package/2006/relationships`">
<Relationship Id=`"rId1`" Type=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument`" Target=`"word/document.xml`"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">
<pkg:xmlData>
<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">
<w:body>
<w:p><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`">docker run hello-world </w:t></w:r><w:r><w:t>=</w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>&gt;  if</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> hello-world image is not found in my local , this command  pull hello- world image pull docker hub and  make it a container and run it</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>"
$ins.InsertXML($xml)
$p2 = $d.Paragraphs.Item(4)
$newEnd = $p2.Range.End - 1
$delStart = $newEnd - $oldLen
$delRng = $d.Range($delStart, $newEnd)
$delRng.Delete()
# --- paragraph 5 ---
$p = $d.Paragraphs.Item(5)
$oldLen = $p.Range.End - $p.Range.Start - 1
$insPoint = $p.Range.Start
$ins = $d.Range($insPoint, $insPoint)
$xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>
<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">
<pkg:part pkg:name=`"/_rels/.rels`" pkg:contentType=`"application/vnd.openxmlformats-package.relationships+xml`" pkg:padding=`"512`">
<pkg:xmlData>
<Relationships xmlns=`"http://schemas.openxmlformats.org/package/2006/relationships`">
<Relationship Id=`"rId1`" Type=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument`" Target=`"word/document.xml`"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">
<pkg:xmlData>
<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">
<w:body>
<w:p><w:r><w:rPr><w:color w:val=`"FF0000`"/></w:rPr><w:t xml:space=`"preserve`">docker build </w:t></w:r><w:r><w:t xml:space=`"preserve`">=&gt; create a </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>image  based</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> on docker file</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>"
$ins.InsertXML($xml)
$p2 = $d.Paragraphs.Item(5)
$newEnd = $p2.Range.End - 1
$delStart = $newEnd - $oldLen
$delRng = $d.Range($delStart, $newEnd)
$delRng.Delete()
Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i=1; $i -le $d.Paragraphs.Count; $i++) { Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text }